$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 49581.223
$ws.Range("I6").Value = 277.6
$ws.Range("J6").Value = 111210.75
$ws.Range("K6").Value = 832.8000000000001
$ws.Range("L6").Value = 333632.25
$ws.Range("M6").Value = -720.8000000000001
$ws.Range("N6").Value = -333856.25

# Row 39
$ws.Range("H39").Value = 402.32144
$ws.Range("I39").Value = 79.46666999999999
$ws.Range("J39").Value = 774.8461
$ws.Range("K39").Value = 238.40001
$ws.Range("L39").Value = 2324.5383
$ws.Range("M39").Value = 57.59999000000002
$ws.Range("N39").Value = -2916.5383

# Row 62
$ws.Range("H62").Value = 1945.375
$ws.Range("I62").Value = 1972.6
$ws.Range("J62").Value = 1900
$ws.Range("K62").Value = 1972.6
$ws.Range("L62").Value = 1900
$ws.Range("M62").Value = -1348.6
$ws.Range("N62").Value = -3148

# Row 65
$ws.Range("H65").Value = 1945.375
$ws.Range("I65").Value = 1972.6
$ws.Range("J65").Value = 1900
$ws.Range("K65").Value = 9863
$ws.Range("L65").Value = 9500
$ws.Range("M65").Value = -6743
$ws.Range("N65").Value = -15740

# Row 92
$ws.Range("H92").Value = 508.5862
$ws.Range("I92").Value = 417.76
$ws.Range("J92").Value = 1076.25
$ws.Range("K92").Value = 417.76
$ws.Range("L92").Value = 1076.25
$ws.Range("M92").Value = 830.24
$ws.Range("N92").Value = -3572.25

# Row 107
$ws.Range("H107").Value = 724.25
$ws.Range("I107").Value = 672.34784
$ws.Range("J107").Value = 963
$ws.Range("K107").Value = 672.34784
$ws.Range("L107").Value = 963
$ws.Range("M107").Value = 1247.65216
$ws.Range("N107").Value = -4803

# Row 132
$ws.Range("H132").Value = 712921.5
$ws.Range("I132").Value = 2525.8333
$ws.Range("J132").Value = 3270346
$ws.Range("K132").Value = 7577.499899999999
$ws.Range("L132").Value = 9811038
$ws.Range("M132").Value = -5047.499899999999
$ws.Range("N132").Value = -9816098

# Row 137
$ws.Range("H137").Value = 1788454.6
$ws.Range("I137").Value = 2858793
$ws.Range("J137").Value = 4557.2856
$ws.Range("K137").Value = 8576379
$ws.Range("L137").Value = 13671.8568
$ws.Range("M137").Value = -8573829
$ws.Range("N137").Value = -18771.8568

# Row 138
$ws.Range("H138").Value = 12504471
$ws.Range("I138").Value = 8598.5
$ws.Range("J138").Value = 16669762
$ws.Range("K138").Value = 25795.5
$ws.Range("L138").Value = 50009286
$ws.Range("M138").Value = -20655.5
$ws.Range("N138").Value = -50019566

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 46340
$ws.Range("I86").Value = 57550
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 57550
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -56427
$ws.Range("N86").Value = -3746

# Row 89
$ws.Range("H89").Value = 46340
$ws.Range("I89").Value = 57550
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 287750
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -282134
$ws.Range("N89").Value = -18732

# Row 94
$ws.Range("H94").Value = 993.5833
$ws.Range("I94").Value = 989
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 989
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -538
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 83339.336
$ws.Range("I17").Value = 50000
$ws.Range("K17").Value = 50000
$ws.Range("M17").Value = -49826

# Row 31
$ws.Range("H31").Value = 283772.03
$ws.Range("I31").Value = 61244.41
$ws.Range("J31").Value = 429270.84
$ws.Range("K31").Value = 61244.41
$ws.Range("L31").Value = 429270.84
$ws.Range("M31").Value = -60949.41
$ws.Range("N31").Value = -429860.84

# Row 34
$ws.Range("H34").Value = 283772.03
$ws.Range("I34").Value = 61244.41
$ws.Range("J34").Value = 429270.84
$ws.Range("K34").Value = 61244.41
$ws.Range("L34").Value = 429270.84
$ws.Range("M34").Value = -61042.41
$ws.Range("N34").Value = -429674.84

# Row 86
$ws.Range("H86").Value = 2527.6
$ws.Range("I86").Value = 1872.2858
$ws.Range("K86").Value = 1872.2858
$ws.Range("M86").Value = -749.2858000000001

# Row 89
$ws.Range("H89").Value = 2527.6
$ws.Range("I89").Value = 1872.2858
$ws.Range("K89").Value = 9361.429
$ws.Range("M89").Value = -3745.429

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372

# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864

# Row 122
$ws.Range("H122").Value = 1108.8462
$ws.Range("I122").Value = 338.8
$ws.Range("J122").Value = 1292.1904
$ws.Range("K122").Value = 3049.2
$ws.Range("L122").Value = 11629.7136
$ws.Range("M122").Value = -599.2000000000003
$ws.Range("N122").Value = -16529.7136

$ws = $wb.Worksheets.Item("GSM")
# Row 69
$ws.Range("H69").Value = 20201
$ws.Range("J69").Value = 20201
$ws.Range("L69").Value = 20201
$ws.Range("N69").Value = -21699

# Row 72
$ws.Range("H72").Value = 20201
$ws.Range("J72").Value = 20201
$ws.Range("L72").Value = 60603
$ws.Range("N72").Value = -68091

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1069.0555
$ws.Range("I16").Value = 856.5
$ws.Range("K16").Value = 856.5
$ws.Range("M16").Value = -686.5

# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 69666.664
$ws.Range("I4").Value = 100000
$ws.Range("J4").Value = 9000
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = -99887
$ws.Range("N4").Value = -9226

# Row 96
$ws.Range("H96").Value = 2154.7273
$ws.Range("I96").Value = 2125.25
$ws.Range("J96").Value = 2233.3333
$ws.Range("K96").Value = 2125.25
$ws.Range("L96").Value = 2233.3333
$ws.Range("M96").Value = -752.25
$ws.Range("N96").Value = -4979.3333

# Row 122
$ws.Range("H122").Value = 958.2941
$ws.Range("I122").Value = 979.06665
$ws.Range("K122").Value = 2937.19995
$ws.Range("M122").Value = -487.1999500000002
